# Inserts a new weekly price row for "Puerro" (Vega Modelo de Temuco) at
# sheet row 186, pushing the existing rows 186-232 down to 187-233 and
# extending the used range to A1:R233.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 186 (shifts 186..232 -> 187..233,
# and copies formatting - e.g. the date number format in column D - from
# the row that used to be there).
$ws.Rows.Item(186).Insert()

# Populate the newly inserted row 186 with the new data point.
$ws.Cells.Item(186, 1).Value2  = 10
$ws.Cells.Item(186, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(186, 3).Value2  = "La Araucanía"
$ws.Cells.Item(186, 4).Value2  = 44798
$ws.Cells.Item(186, 5).Value2  = 9
$ws.Cells.Item(186, 6).Value2  = 100112005
$ws.Cells.Item(186, 7).Value2  = "Puerro"
$ws.Cells.Item(186, 8).Value2  = "Azul de Maquehue"
$ws.Cells.Item(186, 9).Value2  = "Primera"
$ws.Cells.Item(186, 10).Value2 = 30
$ws.Cells.Item(186, 11).Value2 = 16000
$ws.Cells.Item(186, 12).Value2 = 16000
$ws.Cells.Item(186, 13).Value2 = 16000
$ws.Cells.Item(186, 14).Value2 = "$/docena de paquetes"
$ws.Cells.Item(186, 15).Value2 = "Provincia de Cautín"
$ws.Cells.Item(186, 16).Value2 = 1333
$ws.Cells.Item(186, 17).Value2 = 12
$ws.Cells.Item(186, 18).Value2 = "Hortaliza"
